# Applies the scraped-data update for czech-republic_cfl-group-a_2023-2024:
#  - Rows 14/15 swap their match data (F:V)
#  - Rows 29/30/31 rotate their match data (F:V): new29=old30, new30=old31, new31=old29
#  - Rows 76/77 swap their match data (F:V)
#  - Rows 102/103 swap their match data (F:V)
#  - Two new match rows (114, 115) are appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowData($row) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range($c + $row).Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($c in $cols) {
        $ws.Range($c + $row).Value = $data[$c]
    }
}

# --- capture the "before" state of every row that will be rewritten ---
$row14 = Get-RowData 14
$row15 = Get-RowData 15

$row29 = Get-RowData 29
$row30 = Get-RowData 30
$row31 = Get-RowData 31

$row76 = Get-RowData 76
$row77 = Get-RowData 77

$row102 = Get-RowData 102
$row103 = Get-RowData 103

# --- rows 14 / 15 swap ---
Set-RowData 14 $row15
Set-RowData 15 $row14

# --- rows 29 / 30 / 31 rotate (new29=old30, new30=old31, new31=old29) ---
Set-RowData 29 $row30
Set-RowData 30 $row31
Set-RowData 31 $row29

# --- rows 76 / 77 swap ---
Set-RowData 76 $row77
Set-RowData 77 $row76

# --- rows 102 / 103 swap ---
Set-RowData 102 $row103
Set-RowData 103 $row102

# --- append two new rows (114, 115), copying formatting from the last existing row ---
$ws.Range("A113:V113").Copy()
$ws.Range("A114:V115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A114").Value = 113
$ws.Range("B114").Value = "czech-republic"
$ws.Range("C114").Value = "cfl-group-a"
$ws.Range("D114").Value = "2023-2024"
$ws.Range("E114").Value = 45240.75
$ws.Range("F114").Value = "Pisek"
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = "Vltavin"
$ws.Range("I114").Value = 2
$ws.Range("J114").Value = 2.29
$ws.Range("K114").Value = "10/11/2023 11:42"
$ws.Range("L114").Value = 1.83
$ws.Range("M114").Value = "10/11/2023 15:30"
$ws.Range("N114").Value = 3.52
$ws.Range("O114").Value = "10/11/2023 11:42"
$ws.Range("P114").Value = 4.1
$ws.Range("Q114").Value = "10/11/2023 17:33"
$ws.Range("R114").Value = 2.64
$ws.Range("S114").Value = "10/11/2023 11:42"
$ws.Range("T114").Value = 3.4
$ws.Range("U114").Value = "10/11/2023 17:33"
$ws.Range("V114").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/pisek-loko-vltavin/Sd36LVnH/"

$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "czech-republic"
$ws.Range("C115").Value = "cfl-group-a"
$ws.Range("D115").Value = "2023-2024"
$ws.Range("E115").Value = 45240.77083333334
$ws.Range("F115").Value = "Karlovy Vary"
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = "Ceske Budejovice B"
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2.18
$ws.Range("K115").Value = "10/11/2023 11:42"
$ws.Range("L115").Value = 2.64
$ws.Range("M115").Value = "10/11/2023 18:20"
$ws.Range("N115").Value = 3.66
$ws.Range("O115").Value = "10/11/2023 11:42"
$ws.Range("P115").Value = 3.76
$ws.Range("Q115").Value = "10/11/2023 18:21"
$ws.Range("R115").Value = 2.72
$ws.Range("S115").Value = "10/11/2023 11:42"
$ws.Range("T115").Value = 2.27
$ws.Range("U115").Value = "10/11/2023 18:21"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/karlovy-vary-ceske-budejovice/Ei8QQQOd/"

$wb.Save()
